$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.526438280374748
$ws.Range("C2").Value = 0.1308624274661554
$ws.Range("D2").Value = 0.01287503673997747
$ws.Range("F2").Value = 0.5010524358228281
$ws.Range("G2").Value = 0.002381241494065362
$ws.Range("M2").Value = 0.8942543127321727
$ws.Range("N2").Value = 0.9506348847387542
$ws.Range("O2").Value = 1.600791911241117
$ws.Range("B3").Value = 0.4610702097428714
$ws.Range("C3").Value = 0.1191976723034145
$ws.Range("D3").Value = 0.01153341097746363
$ws.Range("F3").Value = 0.491684243643391
$ws.Range("G3").Value = 0.002383930257053685
$ws.Range("M3").Value = 0.7957442636169958
$ws.Range("N3").Value = 0.9649793628531853
$ws.Range("O3").Value = 1.582920526096558
$ws.Range("B4").Value = 0.4208571702748998
$ws.Range("C4").Value = 0.1119834880409201
$ws.Range("D4").Value = 0.01070450451419447
$ws.Range("F4").Value = 0.4863168043954573
$ws.Range("G4").Value = 0.002385669384506046
$ws.Range("M4").Value = 0.7357780737209509
$ws.Range("N4").Value = 0.9742509767348864
$ws.Range("O4").Value = 1.573194046176042
$ws.Range("B5").Value = 0.4044517980485693
$ws.Range("C5").Value = 0.1090308246653251
$ws.Range("D5").Value = 0.01036545166886782
$ws.Range("F5").Value = 0.4842260434287695
$ws.Range("G5").Value = 0.002386400346174901
$ws.Range("M5").Value = 0.7114660064116691
$ws.Range("N5").Value = 0.9781458140287889
$ws.Range("O5").Value = 1.569543088703767
$ws.Range("B6").Value = 0.4017266289700387
$ws.Range("C6").Value = 0.1085397697978863
$ws.Range("D6").Value = 0.01030907647357537
$ws.Range("F6").Value = 0.4838846971401409
$ws.Range("G6").Value = 0.002386523067792691
$ws.Range("M6").Value = 0.7074363620614861
$ws.Range("N6").Value = 0.978799586804378
$ws.Range("O6").Value = 1.568955710344113
$ws.Range("B7").Value = 0.4206359940837672
$ws.Range("C7").Value = 0.1119437189925492
$ws.Range("D7").Value = 0.01069993702296301
$ws.Range("F7").Value = 0.4862882171652529
$ws.Range("G7").Value = 0.002385679152366921
$ws.Range("M7").Value = 0.7354496957839558
$ws.Range("N7").Value = 0.9743030320023927
$ws.Range("O7").Value = 1.573143543212893
$ws.Range("B8").Value = 0.5039159299834353
$ws.Range("C8").Value = 0.1268513516838539
$ws.Range("D8").Value = 0.01241352705842047
$ws.Range("F8").Value = 0.4977423119142372
$ws.Range("G8").Value = 0.002382150312440795
$ws.Range("M8").Value = 0.8601766037815253
$ws.Range("N8").Value = 0.9554844827540165
$ws.Range("O8").Value = 1.594370621745384
$ws.Range("B9").Value = 0.6665789188658664
$ws.Range("C9").Value = 0.1556632325266776
$ws.Range("D9").Value = 0.01573205078423001
$ws.Range("F9").Value = 0.5232675484836307
$ws.Range("G9").Value = 0.002375927017822846
$ws.Range("M9").Value = 1.109171365554417
$ws.Range("N9").Value = 0.9222678589030942
$ws.Range("O9").Value = 1.64592963928547
$ws.Range("B10").Value = 0.7856502145729678
$ws.Range("C10").Value = 0.1765631394149523
$ws.Range("D10").Value = 0.01814352762990268
$ws.Range("F10").Value = 0.5439076343446061
$ws.Range("G10").Value = 0.002371775023741804
$ws.Range("M10").Value = 1.295211018265832
$ws.Range("N10").Value = 0.9001166180872318
$ws.Range("O10").Value = 1.689929551283853
$ws.Range("B11").Value = 0.8397157227526577
$ws.Range("C11").Value = 0.186010594612469
$ws.Range("D11").Value = 0.01923455487974479
$ws.Range("F11").Value = 0.5537112737428345
$ws.Range("G11").Value = 0.002369976470300006
$ws.Range("M11").Value = 1.380612776636013
$ws.Range("N11").Value = 0.8905299617670472
$ws.Range("O11").Value = 1.71128933064179
$ws.Range("B12").Value = 0.8601735069462961
$ws.Range("C12").Value = 0.1895792400424341
$ws.Range("D12").Value = 0.01964681623623221
$ws.Range("F12").Value = 0.5574835266775438
$ws.Range("G12").Value = 0.002369308304025139
$ws.Range("M12").Value = 1.41307062571066
$ws.Range("N12").Value = 0.886970359768263
$ws.Range("O12").Value = 1.719572010614456
$ws.Range("B13").Value = 0.8557682719387003
$ws.Range("C13").Value = 0.1888110685693505
$ws.Range("D13").Value = 0.01955806828996742
$ws.Range("F13").Value = 0.5566684397977752
$ws.Range("G13").Value = 0.002369451632553928
$ws.Range("M13").Value = 1.406074881060292
$ws.Range("N13").Value = 0.8877338382147215
$ws.Range("O13").Value = 1.717779535564603
$ws.Range("B14").Value = 0.8413991173916884
$ws.Range("C14").Value = 0.1863043691544704
$ws.Range("D14").Value = 0.01926848981325691
$ws.Range("F14").Value = 0.55402041883778
$ws.Range("G14").Value = 0.002369921241409458
$ws.Range("M14").Value = 1.383280695875001
$ws.Range("N14").Value = 0.8902356940071279
$ws.Range("O14").Value = 1.711966854412282
$ws.Range("B15").Value = 0.8325955221422987
$ws.Range("C15").Value = 0.1847677762061153
$ws.Range("D15").Value = 0.01909099825924443
$ws.Range("F15").Value = 0.5524062280578192
$ws.Range("G15").Value = 0.002370210569519703
$ws.Range("M15").Value = 1.36933418854187
$ws.Range("N15").Value = 0.8917773601788497
$ws.Range("O15").Value = 1.708431736950303
$ws.Range("B16").Value = 0.782114780723191
$ws.Range("C16").Value = 0.1759444951869114
$ws.Range("D16").Value = 0.01807210389856095
$ws.Range("F16").Value = 0.543275302203952
$ws.Range("G16").Value = 0.002371894373051494
$ws.Range("M16").Value = 1.289645940334196
$ws.Range("N16").Value = 0.900753005646731
$ws.Range("O16").Value = 1.6885607635611
$ws.Range("B17").Value = 0.7511198687880665
$ws.Range("C17").Value = 0.170516133836685
$ws.Range("D17").Value = 0.01744549653226812
$ws.Range("F17").Value = 0.5377800753495023
$ws.Range("G17").Value = 0.002372950390256193
$ws.Range("M17").Value = 1.240962369530763
$ws.Range("N17").Value = 0.9063849194050526
$ws.Range("O17").Value = 1.676715455808932
$ws.Range("B18").Value = 0.7332830148146172
$ws.Range("C18").Value = 0.1673882510199576
$ws.Range("D18").Value = 0.01708452877205957
$ws.Range("F18").Value = 0.5346583532715101
$ws.Range("G18").Value = 0.002373566278302244
$ws.Range("M18").Value = 1.213033036409897
$ws.Range("N18").Value = 0.9096703618656896
$ws.Range("O18").Value = 1.670028795871957
$ws.Range("B19").Value = 0.7272421840094694
$ws.Range("C19").Value = 0.166328244416718
$ws.Range("D19").Value = 0.01696221603710768
$ws.Range("F19").Value = 0.5336080804839298
$ws.Range("G19").Value = 0.002373776268419558
$ws.Range("M19").Value = 1.203588824890176
$ws.Range("N19").Value = 0.910790670844932
$ws.Range("O19").Value = 1.667786496948622
$ws.Range("B20").Value = 0.7544203115394339
$ws.Range("C20").Value = 0.1710945771159231
$ws.Range("D20").Value = 0.01751225810171775
$ws.Range("F20").Value = 0.5383610148592055
$ws.Range("G20").Value = 0.002372837096885726
$ws.Range("M20").Value = 1.246137298429446
$ws.Range("N20").Value = 0.9057806185464372
$ws.Range("O20").Value = 1.677963315454122
$ws.Range("B21").Value = 0.8456201205528373
$ws.Range("C21").Value = 0.1870408912570554
$ws.Range("D21").Value = 0.01935357033553942
$ws.Range("F21").Value = 0.5547965811442879
$ws.Range("G21").Value = 0.002369782956091316
$ws.Range("M21").Value = 1.389972638835104
$ws.Range("N21").Value = 0.8894989191855736
$ws.Range("O21").Value = 1.713668901942867
$ws.Range("B22").Value = 0.9051326976168639
$ws.Range("C22").Value = 0.1974107230043387
$ws.Range("D22").Value = 0.02055179077664349
$ws.Range("F22").Value = 0.5658870071409581
$ws.Range("G22").Value = 0.002367862099375867
$ws.Range("M22").Value = 1.484668914750756
$ws.Range("N22").Value = 0.8792698202784131
$ws.Range("O22").Value = 1.738136918276723
$ws.Range("B23").Value = 0.8733785023368341
$ws.Range("C23").Value = 0.191880996609342
$ws.Range("D23").Value = 0.019912761896272
$ws.Range("F23").Value = 0.559935839400822
$ws.Range("G23").Value = 0.002368880437405971
$ws.Range("M23").Value = 1.434062047371782
$ws.Range("N23").Value = 0.8846915273605269
$ws.Range("O23").Value = 1.724973959155875
$ws.Range("B24").Value = 0.7529282356848626
$ws.Range("C24").Value = 0.1708330848337312
$ws.Range("D24").Value = 0.01748207744591923
$ws.Range("F24").Value = 0.5380982551433675
$ws.Range("G24").Value = 0.002372888289510905
$ws.Range("M24").Value = 1.243797528281561
$ws.Range("N24").Value = 0.9060536747768708
$ws.Range("O24").Value = 1.677398773923272
$ws.Range("B25").Value = 0.6226480380424846
$ws.Range("C25").Value = 0.1479151598713031
$ws.Range("D25").Value = 0.01483889810295835
$ws.Range("F25").Value = 0.5160323409153449
$ws.Range("G25").Value = 0.002377536458338624
$ws.Range("M25").Value = 1.041297463449169
$ws.Range("N25").Value = 0.9308586138364063
$ws.Range("O25").Value = 1.6309115594388
